$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B for the full "Skill Description" name, shifting
# SFIA Level / Keycode / Description one column to the right (C/D/E).
$ws.Columns("B:B").Insert()

# Header
$ws.Range("B1").Value = "Skill Description"

# SkillCode -> Skill Description mapping, applied per data row.
$descriptions = @{
    2  = "Autonomy"
    3  = "Autonomy"
    4  = "Influence"
    5  = "Influence"
    6  = "Influence"
    7  = "Influence"
    8  = "Complexity"
    9  = "Complexity"
    10 = "Complexity"
    11 = "Knowledge"
    12 = "Knowledge"
    13 = "Knowledge"
    15 = "Consultancy"
    16 = "Consultancy"
}

foreach ($row in $descriptions.Keys) {
    $ws.Range("B$row").Value = $descriptions[$row]
}
